$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.334.84"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "1.839.55"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.09"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6257"
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07424"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2893"
$ws.Range("E9").Value = "  -0.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.77"
$ws.Range("E10").Value = "  +1.29%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07718"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").Value = "1.838.98"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.950"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6738"
$ws.Range("E14").Value = "  -1.03%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001020"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.66"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.207"
$ws.Range("E17").Value = "  +0.38%  "
$ws.Range("D18").Value = "29.351.02"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "232.09"
$ws.Range("E19").Value = "  +1.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.28"
$ws.Range("E20").Value = "  -0.52%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.340"
$ws.Range("E22").Value = "  -1.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.10"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.463"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1342"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07272"
$ws.Range("E28").Value = "  +12.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.456"
$ws.Range("E29").Value = "  +4.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.477"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.037"
$ws.Range("E31").Value = "  -1.42%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.032"
$ws.Range("E32").Value = "  -0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.816"
$ws.Range("E33").Value = "  -0.64%  "
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6951"
$ws.Range("E35").Value = "  -0.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.568"
$ws.Range("E36").Value = "  -0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.928"
$ws.Range("E37").Value = "  +4.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01831"
$ws.Range("E38").Value = "  +0.04%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.811"
$ws.Range("D40").Value = "1.229.99"
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9454"
$ws.Range("E41").Value = "  +4.15%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9997"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "1.988.93"
$ws.Range("E43").Value = "  -1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.51"
$ws.Range("E44").Value = "  -0.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.37"
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("E46").Value = "  +5.06%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.708"
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.937"
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.893"
$ws.Range("E49").Value = "  -1.34%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1135"
$ws.Range("E50").Value = "  -3.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3900"
$ws.Range("E51").Value = "  -1.35%  "